$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths (closest achievable values given pixel-snapped COM width granularity)
$ws.Columns.Item(1).ColumnWidth = 15.666666666666666
$ws.Columns.Item(2).ColumnWidth = 14.833333333333334

# Update cell values
$ws.Range("A1").Value = -0.41936842089428694
$ws.Range("B1").Value = 0.41798362865547745
$ws.Range("A2").Value = -0.2650987525301165
$ws.Range("B2").Value = 0.26212040280992532
$ws.Range("A3").Value = -0.15916968011095989
$ws.Range("B3").Value = 0.15829430801895583
$ws.Range("A4").Value = -0.1462943082579109
$ws.Range("B4").Value = 0.14551444810065917
$ws.Range("A5").Value = -0.139514449049849
$ws.Range("B5").Value = 0.13795419929010944
$ws.Range("A6").Value = -0.063115518912184587
$ws.Range("B6").Value = 0.06305247865050978
$ws.Range("A7").Value = -0.043052479793987075
$ws.Range("B7").Value = 0.042941925397210312
$ws.Range("A8").Value = -0.022941926546631741
$ws.Range("B8").Value = 0.022895347146469902
$ws.Range("A9").Value = -0.016895348135440358
$ws.Range("B9").Value = 0.01686824108190077
$ws.Range("A10").Value = -0.01086824207406778
$ws.Range("B10").Value = 0.010870360183055539
$ws.Range("A11").Value = -0.0063703611574865704
$ws.Range("B11").Value = 0.0063712008737297765
$ws.Range("A12").Value = -0.00037120186580397174
$ws.Range("B12").Value = 0.00037057439832377526
$ws.Range("A13").Value = 0.005629424610012812
$ws.Range("B13").Value = -0.005630289060722049
$ws.Range("A14").Value = 0.017630287999029548
$ws.Range("B14").Value = -0.01764531340005604
$ws.Range("A15").Value = 0.023645312411183284
$ws.Range("B15").Value = -0.023678422094265628
$ws.Range("A16").Value = 0.029678421109205821
$ws.Range("B16").Value = -0.02977001084709352
$ws.Range("A17").Value = 0.0068708875109848577
$ws.Range("B17").Value = -0.0068734048155345562
$ws.Range("A18").Value = -0.14288487707075248
$ws.Range("B18").Value = 0.14259194772833439
$ws.Range("A19").Value = -0.082554301402995911
$ws.Range("B19").Value = 0.081633016542719261
$ws.Range("A20").Value = -0.072633017509292408
$ws.Range("B20").Value = 0.072431356763257959
$ws.Range("A21").Value = -0.0090042920858919295
$ws.Range("B21").Value = 0.0089999990265030583
$ws.Range("A22").Value = -0.093950723661045288
$ws.Range("B22").Value = 0.093637358377922197
$ws.Range("A23").Value = -0.084637359348700869
$ws.Range("B23").Value = 0.084127135626684613
$ws.Range("A24").Value = -0.042127136999868142
$ws.Range("B24").Value = 0.04199999861921544
$ws.Range("A25").Value = -0.068838681353863507
$ws.Range("B25").Value = 0.068728756472030028
$ws.Range("A26").Value = -0.062728757448137884
$ws.Range("B26").Value = 0.062591003063797501
$ws.Range("A27").Value = -0.056591004043568205
$ws.Range("B27").Value = 0.056128353099700767
$ws.Range("A28").Value = -0.075288155922429212
$ws.Range("B28").Value = 0.074539996513585116
$ws.Range("A29").Value = -0.062539997591237295
$ws.Range("B29").Value = 0.062170435911552246
$ws.Range("A30").Value = -0.042170437090333301
$ws.Range("B30").Value = 0.042018519006794541
$ws.Range("A31").Value = -0.027018520136174118
$ws.Range("B31").Value = 0.027000423459263345
$ws.Range("A32").Value = -0.0060004246605851463
$ws.Range("B32").Value = 0.0059999989752439475
